$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-43: update Price (D) and/or Volume(1h) (E) columns ---
# Price cells whose new value parses as a plain decimal need NumberFormat
# set to "@" (Text) first so Excel keeps them as text instead of coercing
# them to numbers, matching the original inline-string cell type (these
# "price" values use '.' as a thousands separator, so most of them are not
# valid numbers and naturally stay text).
$ws.Range("D2").Value = "42.269.34"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "2.293.98"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.82"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.03"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.53"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.42"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.958"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.24"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "2.641.72"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "2.284.30"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "42.389.34"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.36"
$ws.Range("E21").Value = "  +33.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.40"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.92"
$ws.Range("E24").Value = "  +7.41%  "
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  +3.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.75"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.40"
$ws.Range("E30").Value = "  +5.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.25"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  +4.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0874"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.65"
$ws.Range("E35").Value = "  -9.11%  "
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.58"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0362"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.70"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.39"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.77"
$ws.Range("E43").Value = "  -2.09%  "

# --- Rows 44 and 45 swapped (Algorand and FirstDigitalUSD exchanged places) ---
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.225"
$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.05%  "

# --- Rows 46-51 ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.96"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.73"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.75"
$ws.Range("E48").Value = "  +4.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.96"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.27"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "1.597.73"
$ws.Range("E51").Value = "  +3.18%  "
